$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits: refine wording of the two indicator cells -----------
# C7 ("Type of decisions" / CMMN column): "events" -> "rules"
$ws.Range("C7").Value = "Stateful conditions, rules, exit and entry criteria"

# B7 ("Type of decisions" / BPMN column): add the missing comma
$ws.Range("B7").Value = "simple, driven by rules or `nevents"

# --- Workbook-level compatibility flag -----------------------------------
try { $wb.CheckCompatibility = $true } catch {}

# --- Add the small (8pt) Calibri font used for phonetic/ruby text -------
# (mirrors turning on "phonetic guide" formatting for the table, which
# registers a dedicated small font without altering the cells' own font)
$titleCell = $ws.Range("A3")
$origSize = $titleCell.Font.Size
$titleCell.Phonetics.Font.Size = 8
$titleCell.Font.Size = $origSize

# --- Sheet view: zoom + active selection ---------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 120
$ws.Range("B5").Select()

# --- Page setup: paper size + orientation --------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
